$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Jerry Lin" to "Jerry L" (row 12, column B)
$ws.Range("B12").Value = "Jerry L"

# Update Saturday's practice attendance answer for rows 10 and 12 to "No"
$ws.Range("D10").Value = "No"
$ws.Range("D12").Value = "No"

# Update the active selection to D11
$ws.Range("D11").Select()
